$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text formatting is preserved (values are stored as text strings, not numbers)
$updates = @{
    "D2" = "330.23"
    "E2" = "-0.43%"
    "D3" = "43.59"
    "E3" = "4.44%"
    "D4" = "5.604"
    "E4" = "-1.77%"
    "D5" = "0.08199"
    "E5" = "-1.70%"
    "D6" = "8.778"
    "E6" = "-0.20%"
    "D7" = "4.400"
    "E7" = "-3.09%"
    "D8" = "1.896"
    "E8" = "-6.96%"
    "E9" = "-4.32%"
    "D10" = "0.9424"
    "E10" = "1.80%"
    "D11" = "0.1193"
    "E11" = "-7.48%"
    "D12" = "0.1924"
    "E12" = "-2.06%"
    "D13" = "0.09930"
    "E13" = "4.45%"
    "D14" = "0.04338"
    "E14" = "10.73%"
    "D15" = "0.1068"
    "E15" = "0.85%"
    "D16" = "0.001287"
    "E16" = "-1.69%"
    "D17" = "0.005999"
    "E17" = "-2.27%"
    "D18" = "3.501"
    "E18" = "1.78%"
    "D20" = "8.718"
    "E20" = "5.73%"
    "D21" = "0.1368"
    "E21" = "-0.26%"
    "E22" = "4.62%"
    "D23" = "0.04381"
    "E23" = "-0.89%"
    "D24" = "0.001238"
    "E24" = "-1.00%"
    "D25" = "0.004307"
    "E25" = "-1.42%"
    "D26" = "0.0001234"
    "E26" = "2.76%"
    "D27" = "0.0004000"
    "E27" = "31.35%"
    "D39" = "0.02802"
    "E39" = "0.05%"
    "D40" = "0.05739"
    "E40" = "3.23%"
    "D41" = "0.007895"
    "E41" = "1.17%"
    "D42" = "0.009767"
    "E42" = "9.25%"
    "D43" = "0.1418"
    "E43" = "-1.09%"
    "D44" = "0.002090"
    "E44" = "-2.40%"
    "D45" = "0.009977"
    "E45" = "-5.28%"
    "D46" = "0.00007305"
    "E46" = "3.69%"
    "D47" = "0.00000000752"
    "E47" = "0.26%"
    "D48" = "0.003501"
    "E48" = "0.14%"
    "D49" = "0.002276"
    "E49" = "-0.08%"
    "D50" = "0.00002106"
    "E50" = "0.26%"
    "D51" = "0.0002006"
    "E51" = "0.26%"
}

foreach ($cellref in $updates.Keys) {
    $range = $ws.Range($cellref)
    $range.NumberFormat = "@"
    $range.Value = $updates[$cellref]
}
